$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1994818652849741
$ws.Range("C2").Value = 0.5544041450777202
$ws.Range("J2").Value = 0.02590673575129534
$ws.Range("P2").Value = 0.1373056994818653
$ws.Range("S2").Value = 0.08290155440414508
$ws.Range("B3").Value = 0.0091324200913242
$ws.Range("C3").Value = 0.0091324200913242
$ws.Range("J3").Value = 0.0319634703196347
$ws.Range("P3").Value = 0.7534246575342466
$ws.Range("S3").Value = 0.1963470319634703
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.2666666666666667
$ws.Range("B6").Value = 0.0641025641025641
$ws.Range("F6").Value = 0.04700854700854701
$ws.Range("J6").Value = 0.2905982905982906
$ws.Range("O6").Value = 0.0170940170940171
$ws.Range("Q6").Value = 0.2051282051282051
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.3205128205128205
$ws.Range("B7").Value = 0.1161616161616162
$ws.Range("D7").Value = 0.04545454545454546
$ws.Range("F7").Value = 0.0303030303030303
$ws.Range("J7").Value = 0.1717171717171717
$ws.Range("O7").Value = 0.04040404040404041
$ws.Range("Q7").Value = 0.2121212121212121
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.2929292929292929
$ws.Range("B8").Value = 0.09937888198757763
$ws.Range("D8").Value = 0.01656314699792961
$ws.Range("F8").Value = 0.06004140786749482
$ws.Range("J8").Value = 0.1118012422360248
$ws.Range("O8").Value = 0.03312629399585922
$ws.Range("Q8").Value = 0.1987577639751553
$ws.Range("R8").Value = 0.07246376811594203
$ws.Range("S8").Value = 0.4078674948240166
$ws.Range("B9").Value = 0.09944751381215469
$ws.Range("D9").Value = 0.02762430939226519
$ws.Range("F9").Value = 0.09392265193370165
$ws.Range("J9").Value = 0.1160220994475138
$ws.Range("O9").Value = 0.03314917127071823
$ws.Range("Q9").Value = 0.1491712707182321
$ws.Range("R9").Value = 0.1104972375690608
$ws.Range("S9").Value = 0.3701657458563536
$ws.Range("B10").Value = 0.134575569358178
$ws.Range("D10").Value = 0.02691511387163561
$ws.Range("F10").Value = 0.06211180124223602
$ws.Range("J10").Value = 0.1111111111111111
$ws.Range("O10").Value = 0.02898550724637681
$ws.Range("Q10").Value = 0.233264320220842
$ws.Range("R10").Value = 0.07453416149068323
$ws.Range("S10").Value = 0.3285024154589372
$ws.Range("G11").Value = 0.152317880794702
$ws.Range("J11").Value = 0.09602649006622517
$ws.Range("K11").Value = 0.2119205298013245
$ws.Range("L11").Value = 0.5264900662251656
$ws.Range("S11").Value = 0.01324503311258278
$ws.Range("G12").Value = 0.7325581395348837
$ws.Range("J12").Value = 0.1976744186046512
$ws.Range("K12").Value = 0.01162790697674419
$ws.Range("L12").Value = 0.04069767441860465
$ws.Range("S12").Value = 0.01744186046511628
$ws.Range("G13").Value = 0.5769230769230769
$ws.Range("J13").Value = 0.3461538461538461
$ws.Range("S13").Value = 0.07692307692307693
$ws.Range("F15").Value = 0.02307692307692308
$ws.Range("H15").Value = 0.1307692307692308
$ws.Range("I15").Value = 0.04230769230769231
$ws.Range("J15").Value = 0.3115384615384615
$ws.Range("K15").Value = 0.06538461538461539
$ws.Range("M15").Value = 0.003846153846153846
$ws.Range("O15").Value = 0.05
$ws.Range("S15").Value = 0.3730769230769231
$ws.Range("F16").Value = 0.0199203187250996
$ws.Range("H16").Value = 0.199203187250996
$ws.Range("I16").Value = 0.05577689243027888
$ws.Range("J16").Value = 0.4581673306772908
$ws.Range("K16").Value = 0.06772908366533864
$ws.Range("M16").Value = 0.0199203187250996
$ws.Range("O16").Value = 0.0398406374501992
$ws.Range("S16").Value = 0.1394422310756972
$ws.Range("F17").Value = 0.02226345083487941
$ws.Range("H17").Value = 0.2133580705009276
$ws.Range("I17").Value = 0.08905380333951762
$ws.Range("J17").Value = 0.4452690166975881
$ws.Range("K17").Value = 0.08348794063079777
$ws.Range("M17").Value = 0.02968460111317254
$ws.Range("N17").Value = 0.001855287569573284
$ws.Range("O17").Value = 0.05009276437847866
$ws.Range("S17").Value = 0.06493506493506493
$ws.Range("F18").Value = 0.0155440414507772
$ws.Range("H18").Value = 0.1347150259067358
$ws.Range("I18").Value = 0.07772020725388601
$ws.Range("J18").Value = 0.4922279792746114
$ws.Range("K18").Value = 0.08290155440414508
$ws.Range("M18").Value = 0.04663212435233161
$ws.Range("O18").Value = 0.05699481865284974
$ws.Range("S18").Value = 0.09326424870466321
$ws.Range("F19").Value = 0.01832669322709163
$ws.Range("H19").Value = 0.2111553784860558
$ws.Range("I19").Value = 0.0749003984063745
$ws.Range("J19").Value = 0.3968127490039841
$ws.Range("K19").Value = 0.1043824701195219
$ws.Range("M19").Value = 0.01673306772908367
$ws.Range("N19").Value = 0.001593625498007968
$ws.Range("O19").Value = 0.07888446215139443
$ws.Range("S19").Value = 0.09721115537848606
